# Weekly data refresh: insert a new price-observation row at row 76
# (shifting all subsequent rows down by one) and populate it with the
# new week's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row before the current row 76; this pushes rows
# 76..163 down to 77..164 and extends the used range/dimension
# automatically (Excel also carries the column-D date style onto the
# freshly inserted row, matching the rest of the column).
$ws.Rows(76).Insert()

# Populate the newly inserted row 76 with the new observation.
$ws.Cells.Item(76, 1).Value  = 9
$ws.Cells.Item(76, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(76, 3).Value  = "Metropolitana"
$ws.Cells.Item(76, 4).Value  = 44494
$ws.Cells.Item(76, 5).Value  = 13
$ws.Cells.Item(76, 6).Value  = 100112043
$ws.Cells.Item(76, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(76, 8).Value  = "Sin especificar"
$ws.Cells.Item(76, 9).Value  = "Primera"
$ws.Cells.Item(76, 10).Value = 85
$ws.Cells.Item(76, 11).Value = 7000
$ws.Cells.Item(76, 12).Value = 7000
$ws.Cells.Item(76, 13).Value = 7000
$ws.Cells.Item(76, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(76, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(76, 16).Value = 117
$ws.Cells.Item(76, 17).Value = 60
$ws.Cells.Item(76, 18).Value = "Hortaliza"
